# Generate Report for handback
# Update the "Latest Handoff Datetime" (D) and "Latest Handback DateTime" (G)
# values for the 980d8046-... entry (row 3) on both the zh-cn and de-de
# localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 19:59:01"
$wsZhCn.Range("G3").Value = "2016-01-08 19:59:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 19:59:11"
$wsDeDe.Range("G3").Value = "2016-01-08 19:59:58"
